# Sprint 1.pptx - remove the "Samenwerkingsovereenkomst" bullet from the
# "Afgelopen Sprint" slide's content placeholder.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(3)
$tr = $shape.TextFrame.TextRange

# Locate the paragraph and delete it (paragraph mark included).
$para = $tr.Paragraphs(4, 1)
$para.Delete()
